$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited range to Text format first so Excel keeps these
# price/volume strings (several of which look like plain numbers,
# e.g. "215.30") as literal text instead of auto-converting them to
# numeric values, matching the original inline-string cell contents.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.111.97"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "1.657.84"
$ws.Range("E3").Value = "  +3.86%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "215.30"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "19.74"
$ws.Range("E10").Value = "  +4.06%  "
$ws.Range("D11").Value = "0.0864"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "1.891.47"
$ws.Range("E12").Value = "  +3.76%  "
$ws.Range("D13").Value = "1.661.54"
$ws.Range("E13").Value = "  +4.04%  "
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  +3.26%  "
$ws.Range("D16").Value = "65.13"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").Value = "27.100.91"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("D18").Value = "238.49"
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "7.90"
$ws.Range("E19").Value = "  +3.09%  "
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("D23").Value = "2.24"
$ws.Range("E23").Value = "  +4.46%  "
$ws.Range("D24").Value = "9.23"
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("D25").Value = "145.64"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.29"
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.518.56"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("E34").Value = "  +4.20%  "
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "0.578"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").Value = "0.891"
$ws.Range("E38").Value = "  +8.65%  "
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "2.27"
$ws.Range("E42").Value = "  +4.46%  "
$ws.Range("E43").Value = "  +9.47%  "
$ws.Range("D44").Value = "1.798.79"
$ws.Range("D45").Value = "0.778"
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("D46").Value = "0.920"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("D50").Value = "0.0506"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "0.0978"
$ws.Range("E51").Value = "  +3.17%  "

# Restore the default (Normal) style so no stray number-format index
# is left attached to the cells - matches the un-styled source cells.
$editRange.Style = "Normal"
